$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" sheet by copying "2022-Q3" (keeps identical formatting) ---
$src = $wb.Worksheets.Item("2022-Q3")
$src.Copy($src)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# --- Step 2: overwrite the 2022-Q4 sheet body with the new quarter data ---
$q4.Range("B2:G13").NumberFormat = "@"
$q4.Range("H2:H13").NumberFormat = "General"

$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 1).Font.Bold = $true
$q4.Cells.Item(2, 1).HorizontalAlignment = -4108
$q4.Cells.Item(2, 1).VerticalAlignment = -4160
$q4.Cells.Item(2, 1).Borders.LineStyle = 1
$q4.Cells.Item(2, 2).Value = "003501"
$q4.Cells.Item(2, 3).Value = "泰达宏利睿智稳健灵活配置混合A"
$q4.Cells.Item(2, 4).Value = "9.84"
$q4.Cells.Item(2, 5).Value = "82.46"
$q4.Cells.Item(2, 6).Value = "2.36"
$q4.Cells.Item(2, 7).Value = "0.2322"
$q4.Cells.Item(2, 8).Value = 6

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 1).Font.Bold = $true
$q4.Cells.Item(3, 1).HorizontalAlignment = -4108
$q4.Cells.Item(3, 1).VerticalAlignment = -4160
$q4.Cells.Item(3, 1).Borders.LineStyle = 1
$q4.Cells.Item(3, 2).Value = "162204"
$q4.Cells.Item(3, 3).Value = "泰达宏利行业精选混合A"
$q4.Cells.Item(3, 4).Value = "8.45"
$q4.Cells.Item(3, 5).Value = "84.56"
$q4.Cells.Item(3, 6).Value = "2.29"
$q4.Cells.Item(3, 7).Value = "0.1935"
$q4.Cells.Item(3, 8).Value = 9

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 1).Font.Bold = $true
$q4.Cells.Item(4, 1).HorizontalAlignment = -4108
$q4.Cells.Item(4, 1).VerticalAlignment = -4160
$q4.Cells.Item(4, 1).Borders.LineStyle = 1
$q4.Cells.Item(4, 2).Value = "013993"
$q4.Cells.Item(4, 3).Value = "中欧光熠一年持有期混合型证券投资基金A"
$q4.Cells.Item(4, 4).Value = "6.13"
$q4.Cells.Item(4, 5).Value = "87.20"
$q4.Cells.Item(4, 6).Value = "2.67"
$q4.Cells.Item(4, 7).Value = "0.1637"
$q4.Cells.Item(4, 8).Value = 10

$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 1).Font.Bold = $true
$q4.Cells.Item(5, 1).HorizontalAlignment = -4108
$q4.Cells.Item(5, 1).VerticalAlignment = -4160
$q4.Cells.Item(5, 1).Borders.LineStyle = 1
$q4.Cells.Item(5, 2).Value = "013280"
$q4.Cells.Item(5, 3).Value = "泰达宏利睿智稳健灵活配置混合C"
$q4.Cells.Item(5, 4).Value = "5.43"
$q4.Cells.Item(5, 5).Value = "82.46"
$q4.Cells.Item(5, 6).Value = "2.36"
$q4.Cells.Item(5, 7).Value = "0.1281"
$q4.Cells.Item(5, 8).Value = 6

$q4.Cells.Item(6, 1).Value = 4
$q4.Cells.Item(6, 1).Font.Bold = $true
$q4.Cells.Item(6, 1).HorizontalAlignment = -4108
$q4.Cells.Item(6, 1).VerticalAlignment = -4160
$q4.Cells.Item(6, 1).Borders.LineStyle = 1
$q4.Cells.Item(6, 2).Value = "010043"
$q4.Cells.Item(6, 3).Value = "天弘安康颐和混合A"
$q4.Cells.Item(6, 4).Value = "14.87"
$q4.Cells.Item(6, 5).Value = "23.17"
$q4.Cells.Item(6, 6).Value = "0.62"
$q4.Cells.Item(6, 7).Value = "0.0922"
$q4.Cells.Item(6, 8).Value = 9

$q4.Cells.Item(7, 1).Value = 5
$q4.Cells.Item(7, 1).Font.Bold = $true
$q4.Cells.Item(7, 1).HorizontalAlignment = -4108
$q4.Cells.Item(7, 1).VerticalAlignment = -4160
$q4.Cells.Item(7, 1).Borders.LineStyle = 1
$q4.Cells.Item(7, 2).Value = "015601"
$q4.Cells.Item(7, 3).Value = "泰达宏利行业精选混合C"
$q4.Cells.Item(7, 4).Value = "3.97"
$q4.Cells.Item(7, 5).Value = "84.56"
$q4.Cells.Item(7, 6).Value = "2.29"
$q4.Cells.Item(7, 7).Value = "0.0909"
$q4.Cells.Item(7, 8).Value = 9

$q4.Cells.Item(8, 1).Value = 6
$q4.Cells.Item(8, 1).Font.Bold = $true
$q4.Cells.Item(8, 1).HorizontalAlignment = -4108
$q4.Cells.Item(8, 1).VerticalAlignment = -4160
$q4.Cells.Item(8, 1).Borders.LineStyle = 1
$q4.Cells.Item(8, 2).Value = "013994"
$q4.Cells.Item(8, 3).Value = "中欧光熠一年持有期混合型证券投资基金C"
$q4.Cells.Item(8, 4).Value = "2.48"
$q4.Cells.Item(8, 5).Value = "87.20"
$q4.Cells.Item(8, 6).Value = "2.67"
$q4.Cells.Item(8, 7).Value = "0.0662"
$q4.Cells.Item(8, 8).Value = 8

$q4.Cells.Item(9, 1).Value = 7
$q4.Cells.Item(9, 1).Font.Bold = $true
$q4.Cells.Item(9, 1).HorizontalAlignment = -4108
$q4.Cells.Item(9, 1).VerticalAlignment = -4160
$q4.Cells.Item(9, 1).Borders.LineStyle = 1
$q4.Cells.Item(9, 2).Value = "006973"
$q4.Cells.Item(9, 3).Value = "太平睿盈混合A"
$q4.Cells.Item(9, 4).Value = "3.17"
$q4.Cells.Item(9, 5).Value = "29.39"
$q4.Cells.Item(9, 6).Value = "1.46"
$q4.Cells.Item(9, 7).Value = "0.0463"
$q4.Cells.Item(9, 8).Value = 8

$q4.Cells.Item(10, 1).Value = 8
$q4.Cells.Item(10, 1).Font.Bold = $true
$q4.Cells.Item(10, 1).HorizontalAlignment = -4108
$q4.Cells.Item(10, 1).VerticalAlignment = -4160
$q4.Cells.Item(10, 1).Borders.LineStyle = 1
$q4.Cells.Item(10, 2).Value = "020034"
$q4.Cells.Item(10, 3).Value = "国泰民安增利债券C"
$q4.Cells.Item(10, 4).Value = "1.06"
$q4.Cells.Item(10, 5).Value = "49.57"
$q4.Cells.Item(10, 6).Value = "2.09"
$q4.Cells.Item(10, 7).Value = "0.0222"
$q4.Cells.Item(10, 8).Value = 10

$q4.Cells.Item(11, 1).Value = 9
$q4.Cells.Item(11, 1).Font.Bold = $true
$q4.Cells.Item(11, 1).HorizontalAlignment = -4108
$q4.Cells.Item(11, 1).VerticalAlignment = -4160
$q4.Cells.Item(11, 1).Borders.LineStyle = 1
$q4.Cells.Item(11, 2).Value = "007669"
$q4.Cells.Item(11, 3).Value = "太平睿盈混合C"
$q4.Cells.Item(11, 4).Value = "0.99"
$q4.Cells.Item(11, 5).Value = "29.39"
$q4.Cells.Item(11, 6).Value = "1.46"
$q4.Cells.Item(11, 7).Value = "0.0145"
$q4.Cells.Item(11, 8).Value = 8

$q4.Cells.Item(12, 1).Value = 10
$q4.Cells.Item(12, 1).Font.Bold = $true
$q4.Cells.Item(12, 1).HorizontalAlignment = -4108
$q4.Cells.Item(12, 1).VerticalAlignment = -4160
$q4.Cells.Item(12, 1).Borders.LineStyle = 1
$q4.Cells.Item(12, 2).Value = "010044"
$q4.Cells.Item(12, 3).Value = "天弘安康颐和混合C"
$q4.Cells.Item(12, 4).Value = "1.17"
$q4.Cells.Item(12, 5).Value = "23.17"
$q4.Cells.Item(12, 6).Value = "0.62"
$q4.Cells.Item(12, 7).Value = "0.0073"
$q4.Cells.Item(12, 8).Value = 9

$q4.Cells.Item(13, 1).Value = 11
$q4.Cells.Item(13, 1).Font.Bold = $true
$q4.Cells.Item(13, 1).HorizontalAlignment = -4108
$q4.Cells.Item(13, 1).VerticalAlignment = -4160
$q4.Cells.Item(13, 1).Borders.LineStyle = 1
$q4.Cells.Item(13, 2).Value = "020033"
$q4.Cells.Item(13, 3).Value = "国泰民安增利债券A"
$q4.Cells.Item(13, 4).Value = "0.20"
$q4.Cells.Item(13, 5).Value = "49.57"
$q4.Cells.Item(13, 6).Value = "2.09"
$q4.Cells.Item(13, 7).Value = "0.0042"
$q4.Cells.Item(13, 8).Value = 10

# --- Step 3: update the "总计" (summary) sheet ---
$zj = $wb.Worksheets.Item("总计")
$zj.Cells.Item(2, 1).Value = 0
$zj.Cells.Item(2, 1).Font.Bold = $true
$zj.Cells.Item(2, 1).HorizontalAlignment = -4108
$zj.Cells.Item(2, 1).VerticalAlignment = -4160
$zj.Cells.Item(2, 1).Borders.LineStyle = 1
$zj.Cells.Item(2, 2).NumberFormat = "@"
$zj.Cells.Item(2, 2).Value = "2022-Q4"
$zj.Cells.Item(2, 3).NumberFormat = "General"
$zj.Cells.Item(2, 3).Value = 12
$zj.Cells.Item(2, 4).NumberFormat = "General"
$zj.Cells.Item(2, 4).Value = 1.06

$zj.Cells.Item(3, 1).Value = 1
$zj.Cells.Item(3, 1).Font.Bold = $true
$zj.Cells.Item(3, 1).HorizontalAlignment = -4108
$zj.Cells.Item(3, 1).VerticalAlignment = -4160
$zj.Cells.Item(3, 1).Borders.LineStyle = 1
$zj.Cells.Item(3, 2).NumberFormat = "@"
$zj.Cells.Item(3, 2).Value = "2022-Q3"
$zj.Cells.Item(3, 3).NumberFormat = "General"
$zj.Cells.Item(3, 3).Value = 10
$zj.Cells.Item(3, 4).NumberFormat = "General"
$zj.Cells.Item(3, 4).Value = 0.71

$zj.Cells.Item(4, 1).Value = 2
$zj.Cells.Item(4, 1).Font.Bold = $true
$zj.Cells.Item(4, 1).HorizontalAlignment = -4108
$zj.Cells.Item(4, 1).VerticalAlignment = -4160
$zj.Cells.Item(4, 1).Borders.LineStyle = 1
$zj.Cells.Item(4, 2).NumberFormat = "@"
$zj.Cells.Item(4, 2).Value = "2022-Q2"
$zj.Cells.Item(4, 3).NumberFormat = "General"
$zj.Cells.Item(4, 3).Value = 10
$zj.Cells.Item(4, 4).NumberFormat = "General"
$zj.Cells.Item(4, 4).Value = 0.29

$zj.Cells.Item(5, 1).Value = 3
$zj.Cells.Item(5, 1).Font.Bold = $true
$zj.Cells.Item(5, 1).HorizontalAlignment = -4108
$zj.Cells.Item(5, 1).VerticalAlignment = -4160
$zj.Cells.Item(5, 1).Borders.LineStyle = 1
$zj.Cells.Item(5, 2).NumberFormat = "@"
$zj.Cells.Item(5, 2).Value = "2022-Q1"
$zj.Cells.Item(5, 3).NumberFormat = "General"
$zj.Cells.Item(5, 3).Value = 14
$zj.Cells.Item(5, 4).NumberFormat = "General"
$zj.Cells.Item(5, 4).Value = 0.38

$zj.Cells.Item(6, 1).Value = 4
$zj.Cells.Item(6, 1).Font.Bold = $true
$zj.Cells.Item(6, 1).HorizontalAlignment = -4108
$zj.Cells.Item(6, 1).VerticalAlignment = -4160
$zj.Cells.Item(6, 1).Borders.LineStyle = 1
$zj.Cells.Item(6, 2).NumberFormat = "@"
$zj.Cells.Item(6, 2).Value = "2021-Q4"
$zj.Cells.Item(6, 3).NumberFormat = "General"
$zj.Cells.Item(6, 3).Value = 24
$zj.Cells.Item(6, 4).NumberFormat = "General"
$zj.Cells.Item(6, 4).Value = 4.86

$zj.Cells.Item(7, 1).Value = 5
$zj.Cells.Item(7, 1).Font.Bold = $true
$zj.Cells.Item(7, 1).HorizontalAlignment = -4108
$zj.Cells.Item(7, 1).VerticalAlignment = -4160
$zj.Cells.Item(7, 1).Borders.LineStyle = 1
$zj.Cells.Item(7, 2).NumberFormat = "@"
$zj.Cells.Item(7, 2).Value = "2021-Q3"
$zj.Cells.Item(7, 3).NumberFormat = "General"
$zj.Cells.Item(7, 3).Value = 17
$zj.Cells.Item(7, 4).NumberFormat = "General"
$zj.Cells.Item(7, 4).Value = 5.19

# --- Step 4: make sure "总计" stays the active/first sheet selection ---
$wb.Worksheets.Item("总计").Activate()
"ok"
